# Applies the resume content edits described by the commit:
#   "Updated resume and assets with new information and formatting changes"
#
# Net textual changes (verified against the OOXML diff):
#   1. "Fullstack" -> "Full-Stack"   (header line, professional summary, and the
#      first Work Experience job title -- 3 occurrences, same rebrand each time)
#   2. "UI/UX Designer | Web Developer" -> "UI/UX Designer | Full-Stack Developer"
#      (second Work Experience job title)
#   3. "December 2023" -> "Present"  (Myteacher Institute end date -> ongoing role)

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

# 1) "Fullstack" -> "Full-Stack" everywhere (header, summary, job title #1).
$d.Content.Find.Execute("Fullstack", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "Full-Stack", $wdReplaceAll)

# 2) Second Work Experience entry's title: "Web Developer" becomes "Full-Stack Developer".
$d.Content.Find.Execute("UI/UX Designer | Web Developer", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "UI/UX Designer | Full-Stack Developer", $wdReplaceAll)

# 3) That same entry is now ongoing: end date "December 2023" becomes "Present".
$d.Content.Find.Execute("December 2023", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Present", $wdReplaceAll)
